# Apply updated odds values to Sheet1, rows 3-9, as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G3" = 2.1;   "H3" = 3.6;   "I3" = 3.2;   "J3" = 2.75;  "L3" = 3.6;
    "U3" = 1.57;  "V3" = 2.25;  "W3" = 10;    "X3" = 12;    "Y3" = 9;
    "Z3" = 19;    "AC3" = 15;   "AD3" = 7;    "AG3" = 12;   "AH3" = 19;
    "AI3" = 11;   "AJ3" = 34;   "AK3" = 23;   "AL3" = 26;   "AM3" = 126;
    "AN3" = 4.33; "AO3" = 11;   "AW3" = 5.5;  "AX3" = 17;   "BB3" = 126;

    "G4" = 1.42;  "I4" = 7.5;   "J4" = 1.95;  "Q4" = 1.99;  "R4" = 1.91;
    "U4" = 2.2;   "V4" = 1.62;  "Y4" = 8.5;   "AG4" = 17;   "AH4" = 41;
    "AI4" = 23;   "AN4" = 3.2;  "AU4" = 10;   "AZ4" = 201;  "BB4" = 451;

    "H5" = 3.6;   "N5" = 10;    "Q5" = 2.05;  "R5" = 1.85;  "AH5" = 26;

    "M6" = 1.11;  "N6" = 6.5;   "O6" = 1.5;   "P6" = 2.63;  "Q6" = 2.5;
    "R6" = 1.5;

    "Q7" = 2.2;   "R7" = 1.65;

    "P9" = 4;     "T9" = 3.28;
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
